# Update plan's ARR and the preloaded choices
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Numeric Inputs")

# Update the ARR (assumed rate of return) values from 7.65% to 7.3%
$ws.Range("C2").Value = 0.073
$ws.Range("C3").Value = 0.073
$ws.Range("C5").Value = 0.073

# Update the preloaded/active selection on the sheet
$ws.Activate()
$ws.Range("C6").Select()
